$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row element: SendingCluster, LigandSymbol, ReceptorSymbol, TargetCluster,
# then columns E..T (16 numeric values) -- total 20 items, index 0..19
$data = @(
  @("ECs", "Icam1", "Spn", "M1", 2, 1, 29.747261, 59.494522, 0.1897261270801772, 0.1356501617082849, 2, 0.6666666666666666, 1.307929333333333, 3.923788, 0.2442423934209346, 0.2442423934209346, 38.90731524822267, 233.443891489336, 0.04633916337254688, 0.03313152016356831),
  @("ECs", "Icam1", "Spn", "M2", 2, 1, 29.747261, 59.494522, 0.1897261270801772, 0.1356501617082849, 3, 1, 4.047117, 12.141351, 0.7557576065790654, 0.7557576065790653, 120.390645696537, 722.343874179222, 0.1433869637076303, 0.1025186415447166),
  @("FAPs", "Icam1", "Spn", "M1", 3, 1, 21.75021766666667, 65.250653, 0.1387214964445808, 0.1487743969271857, 2, 0.6666666666666666, 1.307929333333333, 3.923788, 0.2442423934209346, 0.2442423934209346, 28.44774769261822, 256.029729233564, 0.03388167031055808, 0.03633701478525198),
  @("FAPs", "Icam1", "Spn", "M2", 3, 1, 21.75021766666667, 65.250653, 0.1387214964445808, 0.1487743969271857, 3, 1, 4.047117, 12.141351, 0.7557576065790654, 0.7557576065790653, 88.025675672467, 792.231081052203, 0.1048398261340227, 0.1124373821419337),
  @("M1", "Icam1", "Spn", "M1", 3, 1, 31.00397566666667, 93.011927, 0.1977413728048684, 0.2120713389099788, 2, 0.6666666666666666, 1.307929333333333, 3.923788, 0.2442423934209346, 0.2442423934209346, 40.55100922438622, 364.959083019476, 0.04829682617220237, 0.05179681139135539),
  @("M1", "Icam1", "Spn", "M2", 3, 1, 31.00397566666667, 93.011927, 0.1977413728048684, 0.2120713389099788, 3, 1, 4.047117, 12.141351, 0.7557576065790654, 0.7557576065790653, 125.476716988153, 1129.290452893377, 0.1494445466326661, 0.1602745275186234),
  @("M2", "Icam1", "Spn", "M1", 3, 1, 54.68136733333333, 164.044102, 0.3487542616983074, 0.37402786366769, 2, 0.6666666666666666, 1.307929333333333, 3.923788, 0.2442423934209346, 0.2442423934209346, 71.51936432204178, 643.6742788983761, 0.08518057559294556, 0.09135346062831563),
  @("M2", "Icam1", "Spn", "M2", 3, 1, 54.68136733333333, 164.044102, 0.3487542616983074, 0.37402786366769, 3, 1, 4.047117, 12.141351, 0.7557576065790654, 0.7557576065790653, 221.301891317978, 1991.717021861802, 0.2635736861053618, 0.2826744030393744),
  @("Neutro", "Icam1", "Spn", "M1", 3, 1, 17.571287, 52.713861, 0.1120685440694613, 0.1201899524281925, 2, 0.6666666666666666, 1.307929333333333, 3.923788, 0.2442423934209346, 0.2442423934209346, 22.98200169171867, 206.838015225468, 0.02737188943072471, 0.02935548164621),
  @("Neutro", "Icam1", "Spn", "M2", 3, 1, 17.571287, 52.713861, 0.1120685440694613, 0.1201899524281925, 3, 1, 4.047117, 12.141351, 0.7557576065790654, 0.7557576065790653, 71.11305432957901, 640.017488966211, 0.08469665463873657, 0.0908344707819825),
  @("sCs", "Icam1", "Spn", "M1", 2, 1, 2.0364265, 4.072853, 0.0129881979026049, 0.009286286358668003, 2, 0.6666666666666666, 1.307929333333333, 3.923788, 0.2442423934209346, 0.2442423934209346, 2.663501954527333, 15.981011727164, 0.003172268541956984, 0.002268104806233248),
  @("sCs", "Icam1", "Spn", "M2", 2, 1, 2.0364265, 4.072853, 0.0129881979026049, 0.009286286358668003, 3, 1, 4.047117, 12.141351, 0.7557576065790654, 0.7557576065790653, 8.2416563074005, 49.449937844403, 0.009815929360647919, 0.007018181552434753),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $data[$i]
    $row = 2 + $i

    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]

    for ($col = 5; $col -le 20; $col++) {
        $ws.Cells.Item($row, $col).Value2 = $r[$col - 1]
    }
}

Write-Output "Rows written"
